$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.392.49'
$ws.Range("E2").Value = '  -3.96%  '
$ws.Range("D3").Value = '1.574.48'
$ws.Range("E3").Value = '  -3.21%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").Value = "'289.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.87%  '
$ws.Range("D7").Value = "'0.3677"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.33%  '
$ws.Range("D8").Value = "'49.46"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.39%  '
$ws.Range("D9").Value = "'0.3389"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.59%  '
$ws.Range("D10").Value = "'1.170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.61%  '
$ws.Range("D11").Value = "'0.07621"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.87%  '
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").Value = "'21.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.00%  '
$ws.Range("D14").Value = "'6.064"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.09%  '
$ws.Range("D15").Value = "'6.935"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.66%  '
$ws.Range("D16").Value = "'0.00001139"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.95%  '
$ws.Range("D17").Value = '1.573.24'
$ws.Range("E17").Value = '  -3.55%  '
$ws.Range("D18").Value = "'89.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.20%  '
$ws.Range("D19").Value = "'0.06750"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.69%  '
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").Value = "'6.260"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.90%  '
$ws.Range("D22").Value = "'16.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.70%  '
$ws.Range("D23").Value = "'0.5335"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.69%  '
$ws.Range("D24").Value = "'11.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("D25").Value = '22.419.53'
$ws.Range("E25").Value = '  -3.95%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = "'2.384"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = "'2.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.78%  '
$ws.Range("D28").Value = "'20.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.47%  '
$ws.Range("D29").Value = "'145.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.13%  '
$ws.Range("D30").Value = "'4.978"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.67%  '
$ws.Range("D31").Value = "'125.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.45%  '
$ws.Range("D32").Value = '1.750.32'
$ws.Range("E32").Value = '  -3.38%  '
$ws.Range("D33").Value = "'1.050"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.97%  '
$ws.Range("D34").Value = "'6.299"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.74%  '
$ws.Range("D35").Value = "'1.994"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.78%  '
$ws.Range("D36").Value = "'10.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.23%  '
$ws.Range("D37").Value = "'0.08450"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.70%  '
$ws.Range("D38").Value = "'0.02546"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.81%  '
$ws.Range("D39").Value = "'0.2329"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.55%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = "'0.06568"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.64%  '
$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").Value = "'5.560"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.09%  '
$ws.Range("D42").Value = "'11.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.85%  '
$ws.Range("D43").Value = "'1.249"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.37%  '
$ws.Range("D44").Value = "'0.6389"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.76%  '
$ws.Range("D45").Value = "'14.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.41%  '
$ws.Range("D46").Value = "'0.9993"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("D47").Value = "'0.6009"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.13%  '
$ws.Range("D48").Value = "'3.748"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.62%  '
$ws.Range("D49").Value = "'2.131"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.03%  '
$ws.Range("D50").Value = "'1.261"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.34%  '
$ws.Range("D51").Value = "'123.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.99%  '
